$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2, pushing the existing rows (old 2-10)
# down to (new 3-11). This mirrors the diff, where a brand new job entry was
# inserted above the previous first data row.
$ws.Rows("2:2").Insert()

# Strip any formatting inherited from the row above (the bold header row)
# so the freshly inserted row stays unstyled, like the rest of the data rows.
$ws.Rows("2:2").ClearFormats()

# Populate the new row 2 with the new job entry. Use a leading apostrophe for
# the purely-numeric looking phone number so Excel keeps it as text (not a
# number, which would otherwise drop the content's numeric formatting).
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "'0785972311"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "U092FMBAUP7"
$ws.Range("G2").Value = "0785972311 job no 7896541 pick up zest of india to, office owes 311 £18"

# Clear formats again to drop the "quote prefix" style that gets applied when
# assigning a text value that looks numeric, keeping the row free of any
# explicit cell style (matching the rest of the data rows).
$ws.Rows("2:2").ClearFormats()
